{"js": "const body = context.document.body;\n\n// Each (oldText -> newText) pair corresponds to one answer cell in the table.\n// Pairs are applied in document order so that a value which is produced as a\n// replacement output (e.g. \"18\u00f72=9, 0\") is not re-matched by an earlier rule\n// that also uses it as a search source.\nconst replacements = [\n  [\"62\u00f78=7, 6\", \"45\u00f79=5, 0\"],\n  [\"64\u00f73=21, 1\", \"60\u00f74=15, 0\"],\n  [\"22\u00f73=7, 1\", \"45\u00f75=9, 0\"],\n  [\"38\u00f75=7, 3\", \"58\u00f74=14, 2\"],\n  [\"33\u00f79=3, 6\", \"48\u00f72=24, 0\"],\n  [\"24\u00f79=2, 6\", \"26\u00f79=2, 8\"],\n  [\"94\u00f76=15, 4\", \"73\u00f75=14, 3\"],\n  [\"48\u00f74=12, 0\", \"84\u00f73=28, 0\"],\n  [\"56\u00f79=6, 2\", \"30\u00f79=3, 3\"],\n  [\"42\u00f75=8, 2\", \"43\u00f74=10, 3\"],\n  [\"36\u00f72=18, 0\", \"67\u00f74=16, 3\"],\n  [\"90\u00f79=10, 0\", \"33\u00f72=16, 1\"],\n  [\"18\u00f72=9, 0\", \"78\u00f76=13, 0\"],\n  [\"47\u00f79=5, 2\", \"64\u00f76=10, 4\"],\n  [\"80\u00f79=8, 8\", \"93\u00f76=15, 3\"],\n  [\"83\u00f75=16, 3\", \"88\u00f79=9, 7\"],\n  [\"40\u00f78=5, 0\", \"68\u00f76=11, 2\"],\n  [\"27\u00f76=4, 3\", \"58\u00f73=19, 1\"],\n  [\"12\u00f79=1, 3\", \"76\u00f78=9, 4\"],\n  [\"35\u00f72=17, 1\", \"35\u00f79=3, 8\"],\n  [\"69\u00f74=17, 1\", \"10\u00f79=1, 1\"],\n  [\"33\u00f77=4, 5\", \"86\u00f75=17, 1\"],\n  [\"75\u00f76=12, 3\", \"83\u00f78=10, 3\"],\n  [\"61\u00f74=15, 1\", \"18\u00f72=9, 0\"],\n  [\"82\u00f74=20, 2\", \"90\u00f75=18, 0\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: \"${oldText}\"`);\n  }\n\n  // Replace only the first match (each source string is unique in this\n  // document, but guard against unexpected duplicates anyway).\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# Each (oldText -> newText) pair corresponds to one answer cell in the table.\n# Pairs are applied in document order so that a value which is produced as a\n# replacement output (e.g. \"18\u00f72=9, 0\") is not re-matched by an earlier rule\n# that also uses it as a search source. wdReplaceOne (1) only replaces the\n# first hit found from the start of the story each time.\n$replacements = @(\n  @(\"62\u00f78=7, 6\", \"45\u00f79=5, 0\"),\n  @(\"64\u00f73=21, 1\", \"60\u00f74=15, 0\"),\n  @(\"22\u00f73=7, 1\", \"45\u00f75=9, 0\"),\n  @(\"38\u00f75=7, 3\", \"58\u00f74=14, 2\"),\n  @(\"33\u00f79=3, 6\", \"48\u00f72=24, 0\"),\n  @(\"24\u00f79=2, 6\", \"26\u00f79=2, 8\"),\n  @(\"94\u00f76=15, 4\", \"73\u00f75=14, 3\"),\n  @(\"48\u00f74=12, 0\", \"84\u00f73=28, 0\"),\n  @(\"56\u00f79=6, 2\", \"30\u00f79=3, 3\"),\n  @(\"42\u00f75=8, 2\", \"43\u00f74=10, 3\"),\n  @(\"36\u00f72=18, 0\", \"67\u00f74=16, 3\"),\n  @(\"90\u00f79=10, 0\", \"33\u00f72=16, 1\"),\n  @(\"18\u00f72=9, 0\", \"78\u00f76=13, 0\"),\n  @(\"47\u00f79=5, 2\", \"64\u00f76=10, 4\"),\n  @(\"80\u00f79=8, 8\", \"93\u00f76=15, 3\"),\n  @(\"83\u00f75=16, 3\", \"88\u00f79=9, 7\"),\n  @(\"40\u00f78=5, 0\", \"68\u00f76=11, 2\"),\n  @(\"27\u00f76=4, 3\", \"58\u00f73=19, 1\"),\n  @(\"12\u00f79=1, 3\", \"76\u00f78=9, 4\"),\n  @(\"35\u00f72=17, 1\", \"35\u00f79=3, 8\"),\n  @(\"69\u00f74=17, 1\", \"10\u00f79=1, 1\"),\n  @(\"33\u00f77=4, 5\", \"86\u00f75=17, 1\"),\n  @(\"75\u00f76=12, 3\", \"83\u00f78=10, 3\"),\n  @(\"61\u00f74=15, 1\", \"18\u00f72=9, 0\"),\n  @(\"82\u00f74=20, 2\", \"90\u00f75=18, 0\"),\n)\n\nforeach ($pair in $replacements) {\n  $oldText = $pair[0]\n  $newText = $pair[1]\n\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $find.Text = $oldText\n  $find.Replacement.Text = $newText\n  $found = $find.Execute(\n    $oldText,   # FindText\n    $false,     # MatchCase\n    $false,     # MatchWholeWord\n    $false,     # MatchWildcards\n    $false,     # MatchSoundsLike\n    $false,     # MatchAllWordForms\n    $true,      # Forward\n    1,          # Wrap (wdFindContinue)\n    $false,     # Format\n    $newText,   # ReplaceWith\n    1           # Replace = wdReplaceOne\n  )\n\n  if (-not $found) {\n    throw \"Text not found: $oldText\"\n  }\n}\n"}
